$wb = $excel.ActiveWorkbook

# --- Roster sheet: update pair names ---
$ws1 = $wb.Worksheets.Item("Roster")
$ws1.Range("B8").Value = "Name 44"
$ws1.Range("C8").Value = "Name 46"
$ws1.Range("B9").Value = "Name 49"
$ws1.Range("C9").Value = "Name 36"
$ws1.Range("B10").Value = "Name 54"
$ws1.Range("C10").Value = "Name 23"
$ws1.Range("B11").Value = "Name 19"
$ws1.Range("C11").Value = "Name 72"
$ws1.Range("B12").Value = "Name 58"
$ws1.Range("C12").Value = "Name 83"
$ws1.Range("B16").Value = "Name 52"
$ws1.Range("C16").Value = "Name 46"
$ws1.Range("B17").Value = "Name 61"
$ws1.Range("C17").Value = "Name 36"
$ws1.Range("B18").Value = "Name 36"
$ws1.Range("C18").Value = "Name 69"
$ws1.Range("B19").Value = "Name 33"
$ws1.Range("C19").Value = "Name 38"
$ws1.Range("B20").Value = "Name 48"
$ws1.Range("C20").Value = "Name 18"
# --- By Round sheet: fill in scores for NS/EW columns (J/K) ---
$ws = $wb.Worksheets.Item("By Round")
$ws.Range("J3").Value = 20
$ws.Range("J4").Value = 410
$ws.Range("K5").Value = 70
$ws.Range("J6").Value = 620
$ws.Range("K7").Value = 500
$ws.Range("J8").Value = 300
$ws.Range("J9").Value = 280
$ws.Range("J10").Value = 230
$ws.Range("K11").Value = 520
$ws.Range("J12").Value = "Avg"
$ws.Range("K12").Value = "Avg"
$ws.Range("K13").Value = 220
$ws.Range("K14").Value = 600
$ws.Range("J15").Value = 760
$ws.Range("K16").Value = 730
$ws.Range("J17").Value = 760
$ws.Range("J18").Value = 700
$ws.Range("J19").Value = "Avg"
$ws.Range("K19").Value = "Avg"
$ws.Range("J20").Value = 770
$ws.Range("K21").Value = 90
$ws.Range("K22").Value = 560
$ws.Range("K23").Value = 770
$ws.Range("K24").Value = 40
$ws.Range("K25").Value = 80
$ws.Range("K26").Value = 660
$ws.Range("K27").Value = 370
$ws.Range("K28").Value = 500
$ws.Range("K29").Value = 430
$ws.Range("K30").Value = 30
$ws.Range("K31").Value = 330
$ws.Range("K32").Value = 270
$ws.Range("K33").Value = 780
$ws.Range("J34").Value = 770
$ws.Range("J35").Value = 390
$ws.Range("J36").Value = 460
$ws.Range("K37").Value = 590
$ws.Range("K38").Value = 750
$ws.Range("K39").Value = 640
$ws.Range("J40").Value = 460
$ws.Range("K41").Value = 120
$ws.Range("K42").Value = 490
$ws.Range("J43").Value = 570
$ws.Range("K44").Value = 360
$ws.Range("J45").Value = 530
$ws.Range("J46").Value = "Avg"
$ws.Range("K46").Value = "Avg"
$ws.Range("K47").Value = 280
$ws.Range("J48").Value = 260
$ws.Range("J49").Value = 420
$ws.Range("K50").Value = 700
$ws.Range("K51").Value = 310
$ws.Range("J52").Value = 720
$ws.Range("K53").Value = 270
$ws.Range("J54").Value = "Avg"
$ws.Range("K54").Value = "Avg"
$ws.Range("K55").Value = 700
$ws.Range("K56").Value = 640
$ws.Range("J57").Value = 260
$ws.Range("J58").Value = 50
$ws.Range("J59").Value = "Avg"
$ws.Range("K59").Value = "Avg"
$ws.Range("K60").Value = 190
$ws.Range("J61").Value = "Avg"
$ws.Range("K61").Value = "Avg"
$ws.Range("J62").Value = 170
$ws.Range("K63").Value = 380
$ws.Range("J64").Value = 90
$ws.Range("K65").Value = 160
$ws.Range("J66").Value = 610
$ws.Range("K67").Value = 380
$ws.Range("K68").Value = 620
$ws.Range("J69").Value = "Avg"
$ws.Range("K69").Value = "Avg"
$ws.Range("J70").Value = 610
$ws.Range("J71").Value = 560
$ws.Range("J72").Value = 460
$ws.Range("J73").Value = "Avg"
$ws.Range("K73").Value = "Avg"
$ws.Range("J74").Value = 800
$ws.Range("J75").Value = "Avg"
$ws.Range("K75").Value = "Avg"
$ws.Range("J76").Value = 260
$ws.Range("K77").Value = 780